$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows for the old "ECs" sender (A8:T10) -- only FAPs/MuSCs remain as senders
$ws.Rows("8:10").Delete()

# Update remaining rows (2-7) with refreshed TPM-derived values and re-ordered cluster labels
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Alcam"
$ws.Range("C2").Value = "Chl1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.4771596666666666
$ws.Range("H2").Value = 1.431479
$ws.Range("I2").Value = 0.4973652976730675
$ws.Range("J2").Value = 0.4973652976730676
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.05940166666666667
$ws.Range("N2").Value = 0.178205
$ws.Range("O2").Value = 0.01186516648651917
$ws.Range("P2").Value = 0.01186516648651917
$ws.Range("Q2").Value = 0.02834407946611111
$ws.Range("R2").Value = 0.255096715195
$ws.Range("S2").Value = 0.005901322061508109
$ws.Range("T2").Value = 0.005901322061508111
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Alcam"
$ws.Range("C3").Value = "Chl1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.4771596666666666
$ws.Range("H3").Value = 1.431479
$ws.Range("I3").Value = 0.4973652976730675
$ws.Range("J3").Value = 0.4973652976730676
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.051033
$ws.Range("N3").Value = 0.153099
$ws.Range("O3").Value = 0.01019356989938328
$ws.Range("P3").Value = 0.01019356989938328
$ws.Range("Q3").Value = 0.024350889269
$ws.Range("R3").Value = 0.219158003421
$ws.Range("S3").Value = 0.005069927927357987
$ws.Range("T3").Value = 0.005069927927357988
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Alcam"
$ws.Range("C4").Value = "Chl1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.4771596666666666
$ws.Range("H4").Value = 1.431479
$ws.Range("I4").Value = 0.4973652976730675
$ws.Range("J4").Value = 0.4973652976730676
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4.895956666666667
$ws.Range("N4").Value = 14.68787
$ws.Range("O4").Value = 0.9779412636140976
$ws.Range("P4").Value = 0.9779412636140976
$ws.Range("Q4").Value = 2.336153051081111
$ws.Range("R4").Value = 21.02537745973
$ws.Range("S4").Value = 0.4863940476842015
$ws.Range("T4").Value = 0.4863940476842016
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Alcam"
$ws.Range("C5").Value = "Chl1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.482215
$ws.Range("H5").Value = 1.446645
$ws.Range("I5").Value = 0.5026347023269324
$ws.Range("J5").Value = 0.5026347023269324
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.05940166666666667
$ws.Range("N5").Value = 0.178205
$ws.Range("O5").Value = 0.01186516648651917
$ws.Range("P5").Value = 0.01186516648651917
$ws.Range("Q5").Value = 0.02864437469166667
$ws.Range("R5").Value = 0.257799372225
$ws.Range("S5").Value = 0.005963844425011055
$ws.Range("T5").Value = 0.005963844425011055
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Alcam"
$ws.Range("C6").Value = "Chl1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.482215
$ws.Range("H6").Value = 1.446645
$ws.Range("I6").Value = 0.5026347023269324
$ws.Range("J6").Value = 0.5026347023269324
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.051033
$ws.Range("N6").Value = 0.153099
$ws.Range("O6").Value = 0.01019356989938328
$ws.Range("P6").Value = 0.01019356989938328
$ws.Range("Q6").Value = 0.024608878095
$ws.Range("R6").Value = 0.221479902855
$ws.Range("S6").Value = 0.005123641972025294
$ws.Range("T6").Value = 0.005123641972025294
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Alcam"
$ws.Range("C7").Value = "Chl1"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.482215
$ws.Range("H7").Value = 1.446645
$ws.Range("I7").Value = 0.5026347023269324
$ws.Range("J7").Value = 0.5026347023269324
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.895956666666667
$ws.Range("N7").Value = 14.68787
$ws.Range("O7").Value = 0.9779412636140976
$ws.Range("P7").Value = 0.9779412636140976
$ws.Range("Q7").Value = 2.360903744016667
$ws.Range("R7").Value = 21.24813369615
$ws.Range("S7").Value = 0.4915472159298961
$ws.Range("T7").Value = 0.4915472159298961
